# Task Breakdown List update (Ruchit Jain)
# Adds a 2-hour "Buffer" task (T22) to the SSDMS 46 story, a new Status
# column (H) with DONE / target-date markers, and subtotal rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)
$ws.Activate()

# ---------------------------------------------------------------------
# New column H: width + header
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 23.21875
$ws.Range("H2").Value = "Status "

# ---------------------------------------------------------------------
# SSDMS 46 story (rows 7-29): correct burnt-hours data for a few tasks
# ---------------------------------------------------------------------
$ws.Range("F7").Value = 1
$ws.Range("G7").Formula = "=E7-F7"

$ws.Range("F9").Value = 4
$ws.Range("G9").Formula = "=E9-F9"

$ws.Range("F21").Value = 0.5
$ws.Range("G21").Formula = "=E21-F21"

$ws.Range("F25").Value = 1
$ws.Range("G25").Formula = "=E25-F25"

# ---------------------------------------------------------------------
# New row 29: T22 / Buffer task (2 hours, nothing burnt yet)
# ---------------------------------------------------------------------
$ws.Range("C29").Value = "T22"
$ws.Range("D29").Value = "Buffer"
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 0
$ws.Range("G29").Formula = "=E29-F29"

# ---------------------------------------------------------------------
# Subtotal rows
# ---------------------------------------------------------------------
$ws.Range("G6").Formula = "=SUM(G3:G5)"
$ws.Range("G30").Formula = "=SUM(G7:G29)"
$ws.Range("G47").Formula = "=SUM(G32:G46)"

# ---------------------------------------------------------------------
# Status column content
# ---------------------------------------------------------------------
$ws.Range("H7").Value = "DONE"
$ws.Range("H24").Value = "Complete by 6th July 2017"

# ---------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------

# G5 - plain/no-color font on a white fill
$ws.Range("G5").Interior.ThemeColor = 2
$ws.Range("G5").Font.Color = $ws.Range("G5").Font.Color()

# Subtotal rows: grey fill, top border
$subtotalRanges = @("G6", "G47")
foreach ($r in $subtotalRanges) {
    $ws.Range($r).Interior.ThemeColor = 7
    $ws.Range($r).Borders.Item(8).LineStyle = 1
}
$ws.Range("G30").Interior.ThemeColor = 7

# G28/G29 thin bottom border under the SSDMS 46 story total
$ws.Range("G29").Borders.Item(9).LineStyle = 1

# Status column - DONE block (green tint), centered
$doneRange = $ws.Range("H7:H23")
$doneRange.Interior.ThemeColor = 10
$doneRange.HorizontalAlignment = -4108
$doneRange.VerticalAlignment = -4108
$doneRange.Merge()

# Status column - target-date block (orange tint), vertically centered
$dateRange = $ws.Range("H24:H29")
$dateRange.Interior.ThemeColor = 7
$dateRange.VerticalAlignment = -4108
$dateRange.Merge()

# ---------------------------------------------------------------------
# View state
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("K26").Select()
